$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the first instructions (column D) to be inside the circle (value 3 instead of 2)
$ws.Range("D2").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("D6").Value = 3

# Update the selected cell in the sheet view
$ws.Range("H9").Select()
